# Applies the "Upload new version with timestamp" edit:
#  - Removes the "ZINC OLIVE BABY CREAM 75 GM" line item from the
#    shortage table.
#  - Adds a new line item "معجون سيجنال 50 مل" right after
#    "معجون اسنان سيجنال 120 مل عرض" (keeping the table's alphabetical
#    order), with: stock 15:0, price 35.00, sale price 35.0000, 1:0 txns.
#  - Updates the grand-total cell to reflect the swap.
#  - Bumps the generated-at timestamp string from 8:43 PM to 8:47 PM.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to be stored as TEXT
# (some columns, e.g. L/"حد الطلب" and P/"سعر البيع", carry a numeric
# display format even though the sheet stores literal text in them).
function Set-TextValue($range, [string]$value) {
    $fmt = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = $fmt
}

# Final (post-edit) ordered list of the 16 line items that occupy rows
# 71-86, columns: C=name, H=stock, L=order limit, N=price, P=sale price,
# Q=transactions.
$items = @(
    @{ C = "ZURCAL 40MG 14 GASTRO RESISTANT TAB";    H = "1:0";  L = "1"; N = "96.00"; P = "288.0000"; Q = "3:0" },
    @{ C = "ZURCAL 40MG POWDER FOR I.V. INF. VIAL";   H = "2:0";  L = "1"; N = "84.00"; P = "84.0000";  Q = "1:0" },
    @{ C = "بلسم MINK";                                H = "0:0";  L = "0"; N = "70.00"; P = "70.0000";  Q = "1:0" },
    @{ C = "جهاز محلول ";                              H = "8:0";  L = "0"; N = "20.00"; P = "40.0000";  Q = "2:0" },
    @{ C = "حبايه";                                     H = "0:0";  L = "0"; N = "3.00";  P = "3.0000";   Q = "1:0" },
    @{ C = "سرنجات 3 سم";                              H = "0:0";  L = "0"; N = "2.00";  P = "36.0000";  Q = "18:0" },
    @{ C = "سرنجات 5 سم";                              H = "0:0";  L = "0"; N = "3.00";  P = "18.0000";  Q = "6:0" },
    @{ C = "صابون ديتول اوريجنيال 115 جم";            H = "9:0";  L = "0"; N = "30.00"; P = "30.0000";  Q = "1:0" },
    @{ C = "فرش اسنان اورال بي";                       H = "0:0";  L = "0"; N = "65.00"; P = "65.0000";  Q = "1:0" },
    @{ C = "فرشاة اطفال ريتش ديلي";                    H = "35:0"; L = "0"; N = "15.00"; P = "30.0000";  Q = "2:0" },
    @{ C = "كالونا ";                                   H = "0:0";  L = "0"; N = "15.00"; P = "15.0000";  Q = "1:0" },
    @{ C = "محلول رينجر";                               H = "2:0";  L = "0"; N = "24.00"; P = "24.0000";  Q = "1:0" },
    @{ C = "معجون اسنان سيجنال 120 مل عرض";           H = "0:0";  L = "0"; N = "65.00"; P = "65.0000";  Q = "1:0" },
    @{ C = "معجون سيجنال 50 مل";                       H = "15:0"; L = "0"; N = "35.00"; P = "35.0000";  Q = "1:0" },
    @{ C = "ملح انجليزي";                               H = "5:0";  L = "0"; N = "5.00";  P = "5.0000";   Q = "1:0" },
    @{ C = "مناديل FINE";                               H = "6:0";  L = "0"; N = "35.00"; P = "35.0000";  Q = "1:0" }
)

$row = 71
foreach ($it in $items) {
    $ws.Range("C$row").Value = $it.C
    $ws.Range("H$row").Value = $it.H
    Set-TextValue $ws.Range("L$row") $it.L
    $ws.Range("N$row").Value = $it.N
    Set-TextValue $ws.Range("P$row") $it.P
    $ws.Range("Q$row").Value = $it.Q
    $row++
}

# Grand total (sum of the "price" column) shrinks by 69.00 - 35.00 = 34.00
$ws.Range("P87").Value = 5123.6300000000001

# Footer timestamp: 8:43 PM -> 8:47 PM
$ws.Range("A88").Value = "Saturday, 12 July, 2025 8:47 PM"
